# Adds a new worksheet "7__c0904331-c00d-39b" (Note 7 - Cash at bank and
# in hand) at the end of the workbook, mirroring the layout/formatting
# used by the other "note" sheets (e.g. sheet1.xml / sheet2.xml): a
# bold/bordered header row of 0,1,2 in A1:C1, a "2023 $" / "2022 $"
# column-header row, then the line items.

$wb = $excel.ActiveWorkbook

# A clean, unstyled text cell (used below as a formatting donor so that
# newly-written cells don't pick up a stray implicit number format) and
# the cell that carries the existing header style we want to replicate.
$styleDonor = $wb.Worksheets.Item(1).Range("A3")
$headerDonor = $wb.Worksheets.Item(1).Range("A1:C1")

# Add the new sheet at the very end of the workbook.
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "7__c0904331-c00d-39b"

function Write-TextCell($range, [string]$text) {
    # Forces a literal text write (no numeric/currency/date auto-detection)
    # by switching to a text number format before assigning the value, then
    # pastes plain formatting over the cell afterwards so it doesn't keep
    # an explicit @ number format (matches the plain, unstyled cells used
    # for line-item text/values elsewhere in this workbook).
    $range.NumberFormat = "@"
    $range.Value = $text
    $styleDonor.Copy()
    $range.PasteSpecial(-4122)
}

# Row 1: bold/bordered header cells holding 0, 1, 2 - copy the formatting
# used by every other note sheet's row 1, then set the literal values.
$headerDonor.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2

# Row 2: year headers ("2023 $" / "2022 $").
Write-TextCell $ws.Range("B2") "2023 `$"
Write-TextCell $ws.Range("C2") "2022 `$"

# Column A line-item labels first (matches the shared-string insertion
# order produced by the original edit).
$ws.Range("A3").Value = "Cash at bank and in hand"
$ws.Range("A4").Value = "Term deposit"

# Column B figures.
Write-TextCell $ws.Range("B3") "14502614"
Write-TextCell $ws.Range("B4") "2368905"
Write-TextCell $ws.Range("B5") "16871519"

# Column C figures (C3 and C4 share the same text).
Write-TextCell $ws.Range("C3") "11,203,771 2,321,888"
Write-TextCell $ws.Range("C4") "11,203,771 2,321,888"
Write-TextCell $ws.Range("C5") "13525659"
